# Update the cryptocurrency price/volume table with the latest scraped values.
# (rows 38/39 also swap order: Bittensor now ranks above PEPE)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "63.521.77"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -4.23%  "
$ws.Range("E2").Style = "Normal"
# Row 3
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.089.23"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "  -5.19%  "
$ws.Range("E3").Style = "Normal"
# Row 4
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  +0.02%  "
$ws.Range("E4").Style = "Normal"
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "606.89"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  -1.12%  "
$ws.Range("E5").Style = "Normal"
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "144.08"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -8.56%  "
$ws.Range("E6").Style = "Normal"
# Row 7
$ws.Range("E7").NumberFormat = "@"
$ws.Range("E7").Value = "  +0.06%  "
$ws.Range("E7").Style = "Normal"
# Row 8
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "3.084.48"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -5.33%  "
$ws.Range("E8").Style = "Normal"
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.517"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -5.08%  "
$ws.Range("E9").Style = "Normal"
# Row 10
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -8.08%  "
$ws.Range("E10").Style = "Normal"
# Row 11
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.15"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  -11.05%  "
$ws.Range("E11").Style = "Normal"
# Row 12
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  -6.07%  "
$ws.Range("E12").Style = "Normal"
# Row 13
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  -8.97%  "
$ws.Range("E13").Style = "Normal"
# Row 14
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "34.92"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -10.54%  "
$ws.Range("E14").Style = "Normal"
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "3.597.75"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -5.31%  "
$ws.Range("E15").Style = "Normal"
# Row 16
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  +0.85%  "
$ws.Range("E16").Style = "Normal"
# Row 17
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "63.601.84"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -4.29%  "
$ws.Range("E17").Style = "Normal"
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.087.46"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -5.40%  "
$ws.Range("E18").Style = "Normal"
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "6.76"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -9.08%  "
$ws.Range("E19").Style = "Normal"
# Row 20
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "471.43"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  -6.65%  "
$ws.Range("E20").Style = "Normal"
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "14.48"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  -6.34%  "
$ws.Range("E21").Style = "Normal"
# Row 22
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -7.61%  "
$ws.Range("E22").Style = "Normal"
# Row 23
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  -5.75%  "
$ws.Range("E23").Style = "Normal"
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "13.43"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").NumberFormat = "@"
$ws.Range("E24").Value = "  -8.20%  "
$ws.Range("E24").Style = "Normal"
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "83.05"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  -4.76%  "
$ws.Range("E25").Style = "Normal"
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "1.00"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  +0.06%  "
$ws.Range("E26").Style = "Normal"
# Row 27
$ws.Range("E27").NumberFormat = "@"
$ws.Range("E27").Value = "  -9.31%  "
$ws.Range("E27").Style = "Normal"
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "8.27"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -10.33%  "
$ws.Range("E28").Style = "Normal"
# Row 29
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "2.12"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -11.11%  "
$ws.Range("E29").Style = "Normal"
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "6.63"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -6.52%  "
$ws.Range("E30").Style = "Normal"
# Row 31
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  -0.09%  "
$ws.Range("E31").Style = "Normal"
# Row 32
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "2.71"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -5.44%  "
$ws.Range("E32").Style = "Normal"
# Row 33
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  -17.23%  "
$ws.Range("E33").Style = "Normal"
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "25.93"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -7.25%  "
$ws.Range("E34").Style = "Normal"
# Row 35
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.10"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -4.71%  "
$ws.Range("E35").Style = "Normal"
# Row 36
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -9.29%  "
$ws.Range("E36").Style = "Normal"
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "51.99"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -6.79%  "
$ws.Range("E37").Style = "Normal"
# Row 38
$ws.Range("B38").NumberFormat = "@"
$ws.Range("B38").Value = "Bittensor"
$ws.Range("B38").Style = "Normal"
$ws.Range("C38").NumberFormat = "@"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("C38").Style = "Normal"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "455.49"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -8.09%  "
$ws.Range("E38").Style = "Normal"
# Row 39
$ws.Range("B39").NumberFormat = "@"
$ws.Range("B39").Value = "PEPE"
$ws.Range("B39").Style = "Normal"
$ws.Range("C39").NumberFormat = "@"
$ws.Range("C39").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("C39").Style = "Normal"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0₃0718"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  -8.44%  "
$ws.Range("E39").Style = "Normal"
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.87"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  -15.03%  "
$ws.Range("E40").Style = "Normal"
# Row 41
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.0391"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -7.58%  "
$ws.Range("E41").Style = "Normal"
# Row 42
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -8.15%  "
$ws.Range("E42").Style = "Normal"
# Row 43
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "8.26"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -6.42%  "
$ws.Range("E43").Style = "Normal"
# Row 44
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.808.88"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -6.42%  "
$ws.Range("E44").Style = "Normal"
# Row 45
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  -10.65%  "
$ws.Range("E45").Style = "Normal"
# Row 46
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.21"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -12.18%  "
$ws.Range("E46").Style = "Normal"
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "2.38"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -4.86%  "
$ws.Range("E47").Style = "Normal"
# Row 48
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  +0.03%  "
$ws.Range("E48").Style = "Normal"
# Row 49
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "25.80"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -10.86%  "
$ws.Range("E49").Style = "Normal"
# Row 50
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -5.85%  "
$ws.Range("E50").Style = "Normal"
# Row 51
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "117.72"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -2.66%  "
$ws.Range("E51").Style = "Normal"
